$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) column cells that will hold numeric-looking strings
# so Excel stores them as text instead of re-parsing them as numbers.
$priceCells = @("D2","D3","D5","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D21","D22","D23","D24","D25","D26","D27","D28","D30","D31","D32","D33","D34","D36","D38","D39","D40","D41","D42","D43","D45","D46","D47","D48","D49","D50","D51")
foreach ($cell in $priceCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Update Price (D) and Volume(1h) (E) values per latest crypto data refresh
$ws.Range("D2").Value = '30.463.26'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.878.41'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '247.11'
$ws.Range("E5").Value = '  +5.70%  '
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").Value = '0.4764'
$ws.Range("E7").Value = '  +1.84%  '
$ws.Range("D8").Value = '0.2901'
$ws.Range("E8").Value = '  +1.80%  '
$ws.Range("D9").Value = '0.06522'
$ws.Range("E9").Value = '  +0.84%  '
$ws.Range("D10").Value = '21.86'
$ws.Range("E10").Value = '  +4.15%  '
$ws.Range("D11").Value = '0.07725'
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = '97.25'
$ws.Range("E12").Value = '  +4.13%  '
$ws.Range("D13").Value = '0.7377'
$ws.Range("E13").Value = '  +8.51%  '
$ws.Range("D14").Value = '1.875.97'
$ws.Range("E14").Value = '  +0.79%  '
$ws.Range("D15").Value = '5.133'
$ws.Range("E15").Value = '  +1.90%  '
$ws.Range("D16").Value = '273.03'
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").Value = '30.447.52'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = '13.60'
$ws.Range("E18").Value = '  +2.21%  '
$ws.Range("D19").Value = '0.000007589'
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("E20").Value = '  -0.08%  '
$ws.Range("D21").Value = '2.124.43'
$ws.Range("E21").Value = '  +0.55%  '
$ws.Range("D22").Value = '1.0000'
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").Value = '5.252'
$ws.Range("E23").Value = '  +2.44%  '
$ws.Range("D24").Value = '6.180'
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("D25").Value = '9.335'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").Value = '164.03'
$ws.Range("E26").Value = '  -0.44%  '
$ws.Range("D27").Value = '18.85'
$ws.Range("E27").Value = '  +1.98%  '
$ws.Range("D28").Value = '1.942'
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("E29").Value = '  +0.56%  '
$ws.Range("D30").Value = '0.09947'
$ws.Range("E30").Value = '  +0.12%  '
$ws.Range("D31").Value = '1.520'
$ws.Range("E31").Value = '  +4.93%  '
$ws.Range("D32").Value = '4.311'
$ws.Range("E32").Value = '  +2.28%  '
$ws.Range("D33").Value = '4.067'
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("D34").Value = '0.04786'
$ws.Range("E34").Value = '  +2.61%  '
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("D36").Value = '0.7004'
$ws.Range("E36").Value = '  +2.02%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("D38").Value = '0.01871'
$ws.Range("E38").Value = '  +2.19%  '
$ws.Range("D39").Value = '2.725'
$ws.Range("E39").Value = '  -0.95%  '
$ws.Range("D40").Value = '6.343'
$ws.Range("E40").Value = '  +0.85%  '
$ws.Range("D41").Value = '71.02'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '1.949'
$ws.Range("E42").Value = '  +3.57%  '
$ws.Range("D43").Value = '0.4207'
$ws.Range("E43").Value = '  +4.09%  '
$ws.Range("E44").Value = '  -0.07%  '
$ws.Range("D45").Value = '0.8368'
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").Value = '102.87'
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("D47").Value = '9.258'
$ws.Range("E47").Value = '  +1.89%  '
$ws.Range("D48").Value = '7.087'
$ws.Range("E48").Value = '  +2.31%  '
$ws.Range("D49").Value = '35.65'
$ws.Range("E49").Value = '  +4.86%  '
$ws.Range("D50").Value = '926.68'
$ws.Range("E50").Value = '  -0.42%  '
$ws.Range("D51").Value = '0.05644'
$ws.Range("E51").Value = '  +1.27%  '

# Restore default (Normal) style on the Price cells so only the value changed
foreach ($cell in $priceCells) {
    $ws.Range($cell).Style = "Normal"
}
